$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7
$ws.Range("C2").Value = 0.5838986883763003
$ws.Range("E2").Value = 100

$ws.Range("B3").Value = 0.6899999999999999
$ws.Range("C3").Value = 0.5764360018091361
$ws.Range("E3").Value = 100
